$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: copy/paste operations that reuse EXISTING cell text --------
# (must run before the source cells below get overwritten with new text)

# Row 15 ("Programa:") reuses the existing "01/01/2018" text from row 8,
# preserved as text (not re-interpreted as a date serial) and keeping the
# right column style, by pasting the real source cell.
$ws.Range("B8").Copy()
$ws.Paste($ws.Range("B15"))
$ws.Range("C8").Copy()
$ws.Paste($ws.Range("C15"))

# Row 18 ("Método:") reuses the existing "5840535 - Messias Borges Silva"
# text that currently lives in B13/C13. B18/C18 start out completely empty
# so first pull in the right cell format (matching columns B/C), then paste
# the real shared-string value from B13/C13 (still untouched at this point).
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Paste($ws.Range("B18"))
$ws.Range("C13").Copy()
$ws.Paste($ws.Range("C18"))

# --- Phase 2: plain text overwrites ---------------------------------------

# Row 10: "Objetivos:" body text is replaced by the instructor line
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

# Row 13 gains a label in column A ("Programa resumido:") and a new body
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14 becomes the "Short syllabus" row (text shifted up from row 15)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Sustainability. Environment acts and protocols. Environmental issues. Natural resources and their pollution, Carbon credits, Zero waste concept ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Conventional and renewable sources, Technology and sustainable development, Sustainable urbanization, Industrial Ecology."
$ws.Range("C14").Value = "Sustainability. Environment acts and protocols. Environmental issues. Natural resources and their pollution, Carbon credits, Zero waste concept ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Conventional and renewable sources, Technology and sustainable development, Sustainable urbanization, Industrial Ecology."

$ws.Range("A15").Value = "Programa:"

# Row 16 becomes "Syllabus:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Sustainability- need and concept, challenges,Environment acts and protocols, Global, Regional and Local environmental issues, Natural resources and their pollution, Carbon credits, Zero waste concept  ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Green buildings, Green materials, Energy, Conventional and renewable sources,Technology and sustainable development,Sustainable urbanization, Industrial Ecology."
$ws.Range("C16").Value = "Sustainability- need and concept, challenges,Environment acts and protocols, Global, Regional and Local environmental issues, Natural resources and their pollution, Carbon credits, Zero waste concept  ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Green buildings, Green materials, Energy, Conventional and renewable sources,Technology and sustainable development,Sustainable urbanization, Industrial Ecology."

# Row 17 becomes "Avaliação:" and loses its B/C body text entirely
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

$ws.Range("A18").Value = "Método:"

# Rows 19-21 shift their A-column labels up by one (B/C bodies stay put)
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"

# --- Remove the now-obsolete last row (old row 22, "Bibliografia" body) ---
$ws.Rows.Item(22).Delete()

# --- Row heights -----------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
